$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1097
$ws1.Range("F3").Value = 630
$ws1.Range("F4").Value = 344
$ws1.Range("F5").Value = 4869
$ws1.Range("F6").Value = 507
$ws1.Range("F7").Value = 8896
$ws1.Range("F8").Value = 231
$ws1.Range("F9").Value = 506
$ws1.Range("F10").Value = 72
$ws1.Range("F11").Value = 580
$ws1.Range("F12").Value = 65

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 6
$ws2.Range("F6").Value = 2

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1097
$ws4.Range("F3").Value = 630
$ws4.Range("F4").Value = 344
$ws4.Range("F5").Value = 14
$ws4.Range("F7").Value = 4869
$ws4.Range("F8").Value = 507
$ws4.Range("F10").Value = 8896
$ws4.Range("F11").Value = 231
$ws4.Range("F12").Value = 506
$ws4.Range("F13").Value = 72
$ws4.Range("F14").Value = 5
$ws4.Range("F15").Value = 2
$ws4.Range("F16").Value = 580
$ws4.Range("F17").Value = 65
